$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B55").Value = 3298653
$ws.Range("C55").Value = 35331
$ws.Range("D55").Value = 2687045
$ws.Range("E55").Value = 424263
$ws.Range("F55").Value = 187345
